$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.609176140046208, 4.52093789796921, 10.95258244225863, 22.13368618708452, 3.616665401468534, 18.63086432937966, 13.20513461962098, 15.18509397210819, 19.82095810330227)
    3 = @(3.47692919467207, 4.466589957795994, 10.91052669620021, 22.20200050692867, 3.618799703567573, 18.71246013303271, 12.58541907581515, 14.84167862325303, 19.93336213822134)
    4 = @(3.392475785903279, 4.432699181581706, 10.89011222789852, 22.25280279918085, 3.620178444997235, 18.76889809413721, 12.18873638375351, 14.62933030056536, 20.00901454945488)
    5 = @(3.35727459072193, 4.418766760729616, 10.88315525902816, 22.27571836273161, 3.62075751658453, 18.7934818957683, 12.02319733848939, 14.54254434597655, 20.04150331975277)
    6 = @(3.351383022612468, 4.416446220992124, 10.88208233897054, 22.27965671890874, 3.620854712866178, 18.79765947665695, 11.99548071762194, 14.52812202450586, 20.04699804363744)
    7 = @(3.392004186192727, 4.432511764224956, 10.89001288779983, 22.25310290452247, 3.620186184749131, 18.76922323468041, 12.18651934879963, 14.62816073157963, 20.00944599567703)
    8 = @(3.564270806748548, 4.502312393848196, 10.93696136210141, 22.1553946235533, 3.617387172833915, 18.65767793175987, 12.99491300913654, 15.06706563340304, 19.85833243312551)
    9 = @(3.875042362809257, 4.634657744055766, 11.07170826951699, 22.03463720911965, 3.612437464878757, 18.48960688361884, 14.44528208070203, 15.91079863174562, 19.6150450764377)
    10 = @(4.085476158013946, 4.728608158024298, 11.19623530245869, 21.9898088591356, 3.609126004907723, 18.39751520313037, 15.42149140586387, 16.51386284995552, 19.46917385021145)
    11 = @(4.17708821548979, 4.770530184677735, 11.25827038391169, 21.97906867380579, 3.607689361767705, 18.36253625933962, 15.84506994617534, 16.78330993416702, 19.4100640524269)
    12 = @(4.211171414985968, 4.786279058767311, 11.28251917973613, 21.97639747350581, 3.607155315334328, 18.35029245713132, 16.00244736039725, 16.8845443722353, 19.38873213778465)
    13 = @(4.203858280522621, 4.782892999554241, 11.27726337685845, 21.97691057335082, 3.607269888784462, 18.35288469802022, 15.96868878891262, 16.86277873198434, 19.39327942405875)
    14 = @(4.17990453830257, 4.771828429601207, 11.26025027273759, 21.97882089750241, 3.607645225744074, 18.3615088308053, 15.85807843541726, 16.79165509804086, 19.40828793192505)
    15 = @(4.165152511127319, 4.765034388493846, 11.249927350177, 21.98017301778078, 3.607876428518952, 18.36692206838568, 15.78993079744176, 16.7479829876963, 19.41761831523294)
    16 = @(4.079404816981596, 4.725851276301105, 11.1922881469028, 21.99070569425173, 3.609221292163124, 18.39994095880261, 15.39339031189325, 16.49614707310047, 19.47318346976573)
    17 = @(4.025734828496073, 4.701598312135161, 11.15829660716713, 21.99964544778241, 3.610064151843161, 18.42197347212308, 15.14481687744902, 16.34033085465072, 19.50913398911797)
    18 = @(3.994478742338644, 4.687572296832509, 11.13925381095878, 22.00569551586355, 3.610555511332386, 18.43529633633969, 14.99991757432527, 16.25025454189005, 19.53049341099513)
    19 = @(3.983830124709621, 4.682810485034883, 11.13289402515398, 22.00789967781164, 3.610723007010128, 18.43991868366058, 14.9505287756007, 16.21968115719621, 19.53784212518951)
    20 = @(4.031488221808456, 4.704188047530724, 11.16186258136381, 21.99859974298326, 3.609973748511439, 18.41956070766695, 15.17147790273918, 16.35696556094975, 19.50523637940279)
    21 = @(4.186956962087285, 4.775081851587225, 11.26522702467613, 21.97822185021489, 3.607534709760745, 18.35894846151973, 15.89064993469334, 16.81256825869584, 19.40385095218583)
    22 = @(4.285013086223886, 4.820676315165307, 11.33718720836945, 21.97304250990285, 3.605998799259023, 18.32517833308423, 16.3430270651243, 17.10563219275263, 19.34372402363473)
    23 = @(4.233008528259357, 4.796412071292131, 11.29838383949758, 21.97505991871353, 3.606813240588514, 18.34266497066536, 16.10322039633012, 16.94967836298325, 19.37525051458621)
    24 = @(4.028888358880152, 4.703017484961228, 11.16024884758593, 21.99906967059261, 3.61001459870829, 18.42064947538222, 15.15943064938044, 16.34944654619854, 19.50699633702373)
    25 = @(3.794022898273671, 4.599394744720117, 11.03072672288288, 22.05964306670847, 3.613719144873194, 18.52959699629218, 14.0681777750676, 15.68505545293823, 19.67512644303569)
}

$cols = @("C","D","E","F","G","I","K","M","O")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$row"
        $ws.Range($cellRef).Value = $vals[$i]
    }
}
